$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reference data rows for "issue_level"
# Shared-string table order must match target: "issue_level"(21), "External Issue"(22), "Internal Issue"(23)
$ws.Range("A20").Value = "issue_level"
$ws.Range("C21").Value = "External Issue"
$ws.Range("C20").Value = "Internal Issue"
$ws.Range("B20").Value = 1

$ws.Range("A21").Value = "issue_level"
$ws.Range("B21").Value = 2

# Widen column C to fit the longer text
$ws.Columns.Item(3).AutoFit()

# Update the active selection to match the target workbook state
$ws.Range("C24").Select()
